# Update "想去人数" (interest counts) figures on both the "展览" sheet
# and the combined "全部类型" sheet, per gh-pages data refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates: row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3436
$ws1.Range("F6").Value = 4946
$ws1.Range("F21").Value = 4812
$ws1.Range("F29").Value = 295
$ws1.Range("F33").Value = 105
$ws1.Range("F35").Value = 919
$ws1.Range("F39").Value = 829
$ws1.Range("F40").Value = 911

# Sheet "全部类型" (sheet4) updates: same events, different row offsets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 3436
$ws4.Range("F10").Value = 4946
$ws4.Range("F26").Value = 4812
$ws4.Range("F34").Value = 295
$ws4.Range("F39").Value = 105
$ws4.Range("F40").Value = 919
$ws4.Range("F44").Value = 829
$ws4.Range("F45").Value = 911
